# Roboflow Annotation Report 06/23/2025 - Good Night
# Fill in the weekly progress row (row 42) with the new entry and
# move the active selection to reflect where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D42").Value = "23/6/2025"
$ws.Range("E42").Value = 297
$ws.Range("F42").Value = 629
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 1012
$ws.Range("J42").Value = "N/A"

# Update the view / selection state to match the latest edit position.
$ws.Range("I46").Select()
